# Generate Report for Archive
#
# The localization status changes from "Ready for handoff" to
# "In Translation" for every tracked file, on all three sheets
# (Overview, zh-cn, de-de). Excel's column AutoFit then narrows the
# "Status"-related columns to fit the now-shorter text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: zh-cn (col E) and de-de (col F) status cells for both rows
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

# zh-cn sheet: Status column (C) for both data rows
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

# de-de sheet: Status column (C) for both data rows
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# Re-fit the status columns now that the text is shorter.
$wsOverview.Columns("E:F").ColumnWidth = 12.5
$wsZhCn.Columns("C:C").ColumnWidth = 12.5
$wsDeDe.Columns("C:C").ColumnWidth = 12.5
